$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) for the per-row TP/FP/TN/FN flag columns (H:K) ---
$ws.Range("H1").Value = "TP"
$ws.Range("I1").Value = "FP"
$ws.Range("J1").Value = "TN"
$ws.Range("K1").Value = "FN"

# --- Header row (row 1) for the summary count columns (M:P) ---
$ws.Range("M1").Value = "TP"
$ws.Range("N1").Value = "FP"
$ws.Range("O1").Value = "TN"
$ws.Range("P1").Value = "FN"

# --- Row 2: standalone (non-shared) formulas ---
$ws.Range("H2").Formula = "=IF(AND(A2=-1,F2=-1),""TP"")"
$ws.Range("I2").Formula = "=IF(AND(A2=1,F2=-1),""FP"")"
$ws.Range("J2").Formula = "=IF(AND(A2=1,F2=1),""TN"")"
$ws.Range("K2").Formula = "=IF(AND(A2=-1,F2=1),""FN"")"

# --- Row 2 summary COUNTIF formulas (M2:P2) ---
$ws.Range("M2").Formula = "=COUNTIF(H2:H200,""TP"")"
$ws.Range("N2").Formula = "=COUNTIF(I2:I200,""FP"")"
$ws.Range("O2").Formula = "=COUNTIF(J2:J200,""TN"")"
$ws.Range("P2").Formula = "=COUNTIF(K2:K200,""FN"")"

# --- Rows 3:66 share one formula group per column ---
$ws.Range("H3:H66").Formula = "=IF(AND(A3=-1,F3=-1),""TP"")"
$ws.Range("I3:I66").Formula = "=IF(AND(A3=1,F3=-1),""FP"")"
$ws.Range("J3:J66").Formula = "=IF(AND(A3=1,F3=1),""TN"")"
$ws.Range("K3:K66").Formula = "=IF(AND(A3=-1,F3=1),""FN"")"

# --- Rows 67:68 are a second shared formula group per column ---
$ws.Range("H67:H68").Formula = "=IF(AND(A67=-1,F67=-1),""TP"")"
$ws.Range("I67:I68").Formula = "=IF(AND(A67=1,F67=-1),""FP"")"
$ws.Range("J67:J68").Formula = "=IF(AND(A67=1,F67=1),""TN"")"
$ws.Range("K67:K68").Formula = "=IF(AND(A67=-1,F67=1),""FN"")"

# --- Update the saved selection to match the edited workbook ---
$ws.Range("N3").Select() | Out-Null
